$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "파비뉴스 기사 작성 방식과 기자 채용의 어려움"
$ws.Range("E9").Value = "https://pdsi.pabii.com/pabiinews-how-we-run/#utm_source=rss&utm_medium=rss&utm_campaign=pabiinews-how-we-run"

$ws.Range("D28").Value = "test"
$ws.Range("E28").Value = "https://ropiens.tistory.com/209"

$wb.Save()
